$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Table 5" Heading2 paragraph gets an explicit "no list" numPr
# (ilvl=0 / numId=0) plus a hanging-indent (left=792, hanging=432 twips).
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.ListFormat.RemoveNumbers()
$titlePara.Range.ParagraphFormat.LeftIndent = 39.6
$titlePara.Range.ParagraphFormat.FirstLineIndent = -21.6

# ---------------------------------------------------------------------------
# Change 2: rework the footnote row of "Table 5" -- split the "N, 25th
# Percentile..." sentence into two runs (wrapping "analyses" with
# proofErr gramStart/gramEnd), append a new paragraph with the full
# "Adjusted for pre-specified..." covariate list, and drop the old
# "* P-value < 0.2 ... Benjamini-Hochberg procedure" row entirely.
# ---------------------------------------------------------------------------
$table = $d.Tables(1)

$footnoteRow = $null
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $rowText = $table.Rows($i).Range.Text
    if ($rowText -like "*N, 25th Percentile, and 75th Percentile are from the adjusted analyses*") {
        $footnoteRow = $table.Rows($i)
        break
    }
}

$cell = $footnoteRow.Cells(1)
$targetPara = $cell.Range.Paragraphs(1)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:after="100"/><w:ind w:left="100" w:right="100"/><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t xml:space="preserve">N, 25th Percentile, and 75th Percentile are from the adjusted </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>analyses</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:after="100"/><w:ind w:left="100" w:right="100"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cstheme="majorHAnsi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cstheme="majorHAnsi"/><w:color w:val="000000"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>Adjusted for pre-specified and pre-screened covariates: child sex</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cstheme="majorHAnsi"/><w:color w:val="000000"/><w:sz w:val="14"/><w:szCs w:val="14"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">, child birth order, mother’s age, mother’s height, mother’s education, household food security, number of children &lt; 18 years old in the household, number of people living in the compound, distance (in minutes) to the primary water source, household materials (wall, floor, roof), asset-based household wealth (electricity, wardrobe, table, chair or bench, khat, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cstheme="majorHAnsi"/><w:color w:val="000000"/><w:sz w:val="14"/><w:szCs w:val="14"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>chouki</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cstheme="majorHAnsi"/><w:color w:val="000000"/><w:sz w:val="14"/><w:szCs w:val="14"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>, working radio, working black/white or color television, refrigerator, bicycle, motorcycle, sewing machine, mobile phone, land phone, number of cows, number of goats, number of chickens).</w:t></w:r></w:p>
'@

[void]$targetPara.Range.InsertXML($xml)

$benjaminiRow = $null
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $rowText = $table.Rows($i).Range.Text
    if ($rowText -like "*Benjamini-Hochberg procedure*") {
        $benjaminiRow = $table.Rows($i)
        break
    }
}
if ($benjaminiRow -ne $null) {
    [void]$benjaminiRow.Delete()
}

Write-Output "done"
